# This script updates the cryptocurrency price list on the active worksheet,
# applying fresh "Price" (column D) and "Volume(1h)" (column E) figures, and
# swaps the Monero/Stellar rows (40/41) to reflect their new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells are stored as text, so force text number-format on
# each target cell before writing the value (prevents Excel from silently
# re-interpreting values such as "49.00" or "1.00" as numbers).

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Update Price (D) and Volume(1h) (E) columns for rows with data changes ---
Set-TextValue $ws.Range("D2") '48.209.86'
Set-TextValue $ws.Range("E2") '  +1.83%  '
Set-TextValue $ws.Range("D3") '2.505.82'
Set-TextValue $ws.Range("E3") '  +0.68%  '
Set-TextValue $ws.Range("E4") '  +0.04%  '
Set-TextValue $ws.Range("D5") '320.96'
Set-TextValue $ws.Range("E5") '  -0.07%  '
Set-TextValue $ws.Range("D6") '108.63'
Set-TextValue $ws.Range("E6") '  +0.32%  '
Set-TextValue $ws.Range("D7") '0.528'
Set-TextValue $ws.Range("E7") '  +1.16%  '
Set-TextValue $ws.Range("E8") '  +0.06%  '
Set-TextValue $ws.Range("E9") '  +1.24%  '
Set-TextValue $ws.Range("D10") '39.86'
Set-TextValue $ws.Range("E10") '  +1.75%  '
Set-TextValue $ws.Range("D11") '20.24'
Set-TextValue $ws.Range("E11") '  +10.40%  '
Set-TextValue $ws.Range("D12") '0.0816'
Set-TextValue $ws.Range("E12") '  +0.75%  '
Set-TextValue $ws.Range("E13") '  +0.54%  '
Set-TextValue $ws.Range("E14") '  +0.90%  '
Set-TextValue $ws.Range("D15") '2.902.22'
Set-TextValue $ws.Range("E15") '  +0.85%  '
Set-TextValue $ws.Range("D16") '2.514.77'
Set-TextValue $ws.Range("E16") '  +1.08%  '
Set-TextValue $ws.Range("D17") '0.845'
Set-TextValue $ws.Range("E17") '  +0.26%  '
Set-TextValue $ws.Range("D18") '48.038.28'
Set-TextValue $ws.Range("E18") '  +1.71%  '
Set-TextValue $ws.Range("D19") '13.16'
Set-TextValue $ws.Range("E19") '  +0.28%  '
Set-TextValue $ws.Range("E20") '  +0.03%  '
Set-TextValue $ws.Range("D21") '0.0₃0942'
Set-TextValue $ws.Range("E21") '  +0.67%  '
Set-TextValue $ws.Range("D22") '2.73'
Set-TextValue $ws.Range("E22") '  +2.10%  '
Set-TextValue $ws.Range("D23") '72.12'
Set-TextValue $ws.Range("E23") '  +2.50%  '
Set-TextValue $ws.Range("D24") '276.62'
Set-TextValue $ws.Range("E24") '  +12.85%  '
Set-TextValue $ws.Range("D25") '2.56'
Set-TextValue $ws.Range("E25") '  +0.10%  '
Set-TextValue $ws.Range("E26") '  +0.01%  '
Set-TextValue $ws.Range("D27") '25.87'
Set-TextValue $ws.Range("E27") '  +0.60%  '
Set-TextValue $ws.Range("E28") '  +9.46%  '
Set-TextValue $ws.Range("E29") '  +0.85%  '
Set-TextValue $ws.Range("E30") '  +1.51%  '
Set-TextValue $ws.Range("D31") '35.42'
Set-TextValue $ws.Range("E31") '  +2.60%  '
Set-TextValue $ws.Range("D32") '49.00'
Set-TextValue $ws.Range("E32") '  -1.61%  '
Set-TextValue $ws.Range("D33") '19.37'
Set-TextValue $ws.Range("E33") '  -5.46%  '
Set-TextValue $ws.Range("E34") '  +0.22%  '
Set-TextValue $ws.Range("E35") '  -0.07%  '
Set-TextValue $ws.Range("D36") '0.0784'
Set-TextValue $ws.Range("E36") '  -0.10%  '
Set-TextValue $ws.Range("E37") '  -0.07%  '
Set-TextValue $ws.Range("D38") '4.62'
Set-TextValue $ws.Range("E38") '  -2.19%  '
Set-TextValue $ws.Range("E39") '  +1.24%  '
Set-TextValue $ws.Range("D42") '2.21'
Set-TextValue $ws.Range("E42") '  -0.65%  '
Set-TextValue $ws.Range("E43") '  -6.43%  '
Set-TextValue $ws.Range("D44") '0.0306'
Set-TextValue $ws.Range("E44") '  +3.38%  '
Set-TextValue $ws.Range("D45") '2.008.19'
Set-TextValue $ws.Range("E45") '  +0.53%  '
Set-TextValue $ws.Range("D46") '3.14'
Set-TextValue $ws.Range("E46") '  +3.66%  '
Set-TextValue $ws.Range("D47") '1.86'
Set-TextValue $ws.Range("E47") '  +4.48%  '
Set-TextValue $ws.Range("E48") '  -0.34%  '
Set-TextValue $ws.Range("D49") '9.03'
Set-TextValue $ws.Range("E49") '  -1.35%  '
Set-TextValue $ws.Range("E50") '  +1.90%  '
Set-TextValue $ws.Range("D51") '79.84'
Set-TextValue $ws.Range("E51") '  +2.88%  '

# --- Rows 40/41: Monero moves up to rank 40 (row 40), Stellar moves down to
#     rank 41 (row 41); Monero price/volume also updated ---
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D40") '123.15'
Set-TextValue $ws.Range("E40") '  +4.99%  '

$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D41") '0.112'
Set-TextValue $ws.Range("E41") '  +0.49%  '
